# Update Leve market price/profit columns (H:N) across multiple craft job sheets.
# Values refreshed from the scheduled market-data runner; see commit message.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 33: Glazed and Confused / Clear Glass Lens
$ws.Range("H33").Value = 686.6667
$ws.Range("I33").Value = 686.6667
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 686.6667
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = -457.6667
$ws.Range("N33").ClearContents()

# Row 40: Stuck in the Moment / Horn Glue
$ws.Range("H40").Value = 2900
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 2900
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 2900
$ws.Range("N40").Value = -3250
$ws.Range("M40").ClearContents()

# Row 41: The Write Stuff / Enchanted Mythril Ink
$ws.Range("H41").Value = 579.8333
$ws.Range("I41").Value = 740
$ws.Range("J41").Value = 259.5
$ws.Range("K41").Value = 740
$ws.Range("L41").Value = 259.5
$ws.Range("M41").Value = -300

# Row 44: Alive and Unwell / Budding Oak Wand
$ws.Range("H44").Value = 39998.668
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 39998.668
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 39998.668
$ws.Range("N44").Value = -40922.668

# Row 80: Cleansing the Wicked Humours / Hallowed Water
$ws.Range("H80").Value = 637.55554
$ws.Range("I80").Value = 518.5
$ws.Range("J80").Value = 732.8
$ws.Range("K80").Value = 1555.5
$ws.Range("L80").Value = 2198.4
$ws.Range("M80").Value = -557.5
$ws.Range("N80").Value = -4194.4

# Row 83: Washing Away the Sins (L) / Hallowed Water
$ws.Range("H83").Value = 637.55554
$ws.Range("I83").Value = 518.5
$ws.Range("J83").Value = 732.8
$ws.Range("K83").Value = 4666.5
$ws.Range("L83").Value = 6595.2
$ws.Range("M83").Value = 325.5
$ws.Range("N83").Value = -16579.2

# Row 103: Let Loose the Juice / Persimmon Tannin
$ws.Range("H103").Value = 783
$ws.Range("I103").Value = 700
$ws.Range("J103").Value = 824.5
$ws.Range("K103").Value = 2100
$ws.Range("L103").Value = 2473.5
$ws.Range("M103").Value = -1514
$ws.Range("N103").Value = -3645.5

# Row 108: Keeping Magic Alive / Smilodonskin Grimoire
$ws.Range("H108").Value = 0
$ws.Range("I108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("K108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
# Row 63: Rivets Run through It / Mythrite Rivets
$ws.Range("H63").Value = 3577.2222
$ws.Range("I63").Value = 3577.2222
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 3577.2222
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = -2891.2222

# Row 66: A Riveting Revival (L) / Mythrite Rivets
$ws.Range("H66").Value = 3577.2222
$ws.Range("I66").Value = 3577.2222
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 17886.111
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = -14454.111

# Row 88: The Mast Chance / Adamantite Rivets
$ws.Range("H88").Value = 2900.6667
$ws.Range("I88").Value = 2995
$ws.Range("J88").Value = 2853.5
$ws.Range("K88").Value = 2995
$ws.Range("L88").Value = 2853.5
$ws.Range("M88").Value = -2589
$ws.Range("N88").Value = -3665.5

# Row 91: The Rose and the Riveter (L) / Adamantite Rivets
$ws.Range("H91").Value = 2900.6667
$ws.Range("I91").Value = 2995
$ws.Range("J91").Value = 2853.5
$ws.Range("K91").Value = 2995
$ws.Range("L91").Value = 2853.5
$ws.Range("M91").Value = -1591
$ws.Range("N91").Value = -5661.5

# Row 94: Setting the Stage / High Steel Helm of Maiming
$ws.Range("H94").Value = 0
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
# Row 86: Through Thick and Thin / Adamantite Nugget
$ws.Range("H86").Value = 3303
$ws.Range("I86").Value = 3303
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 3303
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -2180
$ws.Range("N86").ClearContents()

# Row 89: Piercing Eyes Deserve Piercing Shafts (L) / Adamantite Nugget
$ws.Range("H89").Value = 3303
$ws.Range("I89").Value = 3303
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 16515
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -10899
$ws.Range("N89").ClearContents()

# Row 122: To Delight a Dancer / High Durium Tathlums
$ws.Range("H122").Value = 55000
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 55000
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 55000
$ws.Range("N122").Value = -64800

$ws = $wb.Worksheets.Item("CRP")
# Row 2: In with the New / Bone Harpoon
$ws.Range("H2").Value = 66.333336
$ws.Range("I2").Value = 49.5
$ws.Range("J2").Value = 100
$ws.Range("K2").Value = 49.5
$ws.Range("L2").Value = 100
$ws.Range("M2").Value = 63.5

# Row 22: Driving Up the Wall / Elm Lumber
$ws.Range("H22").Value = 1000
$ws.Range("I22").Value = 1000
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 1000
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -650

# Row 58: You Do the Heavy Lifting / Mahogany Lumber
$ws.Range("H58").Value = 7989.433
$ws.Range("I58").Value = 742.55554
$ws.Range("J58").Value = 11095.238
$ws.Range("K58").Value = 742.55554
$ws.Range("L58").Value = 11095.238
$ws.Range("M58").Value = -539.55554
$ws.Range("N58").Value = -11501.238

# Row 134: Wood You Be Quiet / Ceiba Lumber
$ws.Range("H134").Value = 6835.4
$ws.Range("I134").Value = 4294.25
$ws.Range("J134").Value = 17000
$ws.Range("K134").Value = 12882.75
$ws.Range("L134").Value = 51000
$ws.Range("M134").Value = -10347.75

# Row 136: Turali Quality / Dark Mahogany Lumber
$ws.Range("H136").Value = 7989.433
$ws.Range("I136").Value = 742.55554
$ws.Range("J136").Value = 11095.238
$ws.Range("K136").Value = 2227.66662
$ws.Range("L136").Value = 33285.714
$ws.Range("M136").Value = 322.33338
$ws.Range("N136").Value = -38385.714

$ws = $wb.Worksheets.Item("CUL")
# Row 12: Butter Me Up / Kukuru Butter
$ws.Range("H12").Value = 29.916666
$ws.Range("I12").Value = 36.8
$ws.Range("J12").Value = 25
$ws.Range("K12").Value = 110.4
$ws.Range("L12").Value = 75
$ws.Range("M12").Value = 62.60000000000001
$ws.Range("N12").Value = -421

# Row 40: True Grits / Cornmeal
$ws.Range("H40").Value = 95.833336
$ws.Range("I40").Value = 65
$ws.Range("J40").Value = 126.666664
$ws.Range("K40").Value = 260
$ws.Range("L40").Value = 506.666656
$ws.Range("M40").Value = -191
$ws.Range("N40").Value = -644.666656

# Row 99: A Shorlonging for the Familiar / Shorlog
$ws.Range("H99").Value = 1320.6666
$ws.Range("I99").Value = 1320.6666
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 3961.9998
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -1715.9998

# Row 103: West Meats East / Nomad Meat Pie
$ws.Range("H103").Value = 2085.6
$ws.Range("I103").Value = 150
$ws.Range("J103").Value = 2569.5
$ws.Range("K103").Value = 450
$ws.Range("L103").Value = 7708.5
$ws.Range("M103").Value = 429

$ws = $wb.Worksheets.Item("GSM")
# Row 70: Sky Is the Limit / Mythrite Ingot
$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("M70").ClearContents()

# Row 73: Hulls of Broken Dreams (L) / Mythrite Ingot
$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("M73").ClearContents()

# Row 102: Put the Metal to the Peddle / Durium Ingot
$ws.Range("H102").Value = 2950.3076
$ws.Range("I102").Value = 2535.4
$ws.Range("J102").Value = 4333.3335
$ws.Range("K102").Value = 2535.4
$ws.Range("L102").Value = 4333.3335
$ws.Range("M102").Value = -913.4000000000001

$ws = $wb.Worksheets.Item("LTW")
# Row 35: No Risk, No Reward / Toadskin Cesti
$ws.Range("H35").Value = 5497.5
$ws.Range("I35").Value = 5497.5
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 5497.5
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = -5161.5

# Row 54: Not So Alike in Dignity / Boarskin Jerkin
$ws.Range("H54").Value = 2076
$ws.Range("I54").Value = 2076
$ws.Range("J54").Value = 0
$ws.Range("K54").Value = 2076
$ws.Range("L54").Value = 0
$ws.Range("M54").Value = -1432
$ws.Range("N54").ClearContents()

# Row 100: Tiger in the Sack / Tiger Leather
$ws.Range("H100").Value = 6314
$ws.Range("I100").Value = 5221.5
$ws.Range("J100").Value = 8499
$ws.Range("K100").Value = 5221.5
$ws.Range("L100").Value = 8499
$ws.Range("M100").Value = -4680.5
$ws.Range("N100").Value = -9581

$ws = $wb.Worksheets.Item("WVR")
# Row 39: By the Short Hairs / Velveteen Robe
$ws.Range("H39").Value = 0
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("N39").ClearContents()

# Row 46: Crunching the Numbers / Linen Hat
$ws.Range("H46").Value = 0
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("M46").ClearContents()

# Row 134: Cloth for Canvas / Mountain Linen
$ws.Range("H134").Value = 0
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("M134").ClearContents()
